# Adding search test cases
# The "Test Cases" sheet's Runmode column (C2:C11) is switched from "Y" to "N"
# for the existing rows, so these cases are skipped while the new search
# test cases are authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("C2:C11").Value = "N"
